# Auto-committed on 2022/03/14 週一
#
# Reproduces the authoring edit captured in the diff:
#  - delete DBD!row 22 (the "RegisteredAddress" / "雙掛號" row), which shifts
#    rows 23-41 up by one and drops the two now-unused shared strings
#  - make DBD (sheet 1) the active/selected sheet, with the cell cursor
#    parked on G22 (was C17 on DBD before, while DBS used to be the active
#    tab)

$wb = $excel.ActiveWorkbook

$wsDBD = $wb.Worksheets.Item("DBD")

# Delete the row that held "RegisteredAddress" / "雙掛號" (old row 22).
# Everything below it shifts up by one row.
$wsDBD.Rows.Item(22).Delete()

# Switch the active tab from DBS to DBD and move the selection to G22.
$wsDBD.Activate()
$wsDBD.Range("G22").Select()
